# Auto-generated edit script: updates fetched market-price columns (H-N)
# in the Kujata_Profits leve-profit tables, per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (39 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3449.9111  # H64: 3440.1086 -> 3449.9111
$ws.Cells.Item(64, 9).Value = 3344.2273  # I64: 3329.2173 -> 3344.2273
$ws.Cells.Item(64, 11).Value = 3344.2273  # K64: 3329.2173 -> 3344.2273
$ws.Cells.Item(64, 13).Value = -3096.2273  # M64: -3081.2173 -> -3096.2273
$ws.Cells.Item(67, 8).Value = 3449.9111  # H67: 3440.1086 -> 3449.9111
$ws.Cells.Item(67, 9).Value = 3344.2273  # I67: 3329.2173 -> 3344.2273
$ws.Cells.Item(67, 11).Value = 3344.2273  # K67: 3329.2173 -> 3344.2273
$ws.Cells.Item(67, 13).Value = -2486.2273  # M67: -2471.2173 -> -2486.2273
$ws.Cells.Item(112, 8).Value = 2477.5  # H112: 2426.5518 -> 2477.5
$ws.Cells.Item(112, 10).Value = 2583.5  # J112: 2524.8518 -> 2583.5
$ws.Cells.Item(112, 12).Value = 7750.5  # L112: 7574.555399999999 -> 7750.5
$ws.Cells.Item(112, 14).Value = -9966.5  # N112: -9790.555399999999 -> -9966.5
$ws.Cells.Item(116, 8).Value = 2819.28  # H116: 2754.8518 -> 2819.28
$ws.Cells.Item(116, 9).Value = 2381.2144  # I116: 2327.25 -> 2381.2144
$ws.Cells.Item(116, 11).Value = 2381.2144  # K116: 2327.25 -> 2381.2144
$ws.Cells.Item(116, 13).Value = 1060.7856  # M116: 1114.75 -> 1060.7856
$ws.Cells.Item(121, 8).Value = 1114.1666  # H121: 1180.9 -> 1114.1666
$ws.Cells.Item(121, 10).Value = 1114.1666  # J121: 1180.9 -> 1114.1666
$ws.Cells.Item(121, 12).Value = 3342.4998  # L121: 3542.7 -> 3342.4998
$ws.Cells.Item(121, 14).Value = -6836.4998  # N121: -7036.700000000001 -> -6836.4998
$ws.Cells.Item(132, 8).Value = 6947154.5  # H132: 7094945 -> 6947154.5
$ws.Cells.Item(132, 9).Value = 9011807  # I132: 9262107 -> 9011807
$ws.Cells.Item(132, 11).Value = 27035421  # K132: 27786321 -> 27035421
$ws.Cells.Item(132, 13).Value = -27032891  # M132: -27783791 -> -27032891
$ws.Cells.Item(137, 8).Value = 1791.2222  # H137: 1791.6842 -> 1791.2222
$ws.Cells.Item(137, 10).Value = 2925  # J137: 2700 -> 2925
$ws.Cells.Item(137, 12).Value = 8775  # L137: 8100 -> 8775
$ws.Cells.Item(137, 14).Value = -13875  # N137: -13200 -> -13875
$ws.Cells.Item(138, 8).Value = 1776.52  # H138: 1795.67 -> 1776.52
$ws.Cells.Item(138, 9).Value = 809.0625  # I138: 833.26666 -> 809.0625
$ws.Cells.Item(138, 10).Value = 1960.7976  # J138: 1965.5059 -> 1960.7976
$ws.Cells.Item(138, 11).Value = 2427.1875  # K138: 2499.79998 -> 2427.1875
$ws.Cells.Item(138, 12).Value = 5882.392800000001  # L138: 5896.5177 -> 5882.392800000001
$ws.Cells.Item(138, 13).Value = 2712.8125  # M138: 2640.20002 -> 2712.8125
$ws.Cells.Item(138, 14).Value = -16162.3928  # N138: -16176.5177 -> -16162.3928
$ws.Cells.Item(141, 8).Value = 1486.5385  # H141: 1501.25 -> 1486.5385
$ws.Cells.Item(141, 9).Value = 1486.5385  # I141: 1501.25 -> 1486.5385
$ws.Cells.Item(141, 11).Value = 4459.6155  # K141: 4503.75 -> 4459.6155
$ws.Cells.Item(141, 13).Value = 720.3845000000001  # M141: 676.25 -> 720.3845000000001

# --- Sheet: ARM (42 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2663.1538  # H32: 2808.8276 -> 2663.1538
$ws.Cells.Item(32, 9).Value = 2842.074  # I32: 3101.4695 -> 2842.074
$ws.Cells.Item(32, 10).Value = 1784.8182  # J32: 1215.5555 -> 1784.8182
$ws.Cells.Item(32, 11).Value = 2842.074  # K32: 3101.4695 -> 2842.074
$ws.Cells.Item(32, 12).Value = 1784.8182  # L32: 1215.5555 -> 1784.8182
$ws.Cells.Item(32, 13).Value = -2555.074  # M32: -2814.4695 -> -2555.074
$ws.Cells.Item(32, 14).Value = -2358.8182  # N32: -1789.5555 -> -2358.8182
$ws.Cells.Item(61, 8).Value = 1271.1111  # H61: 1463.9 -> 1271.1111
$ws.Cells.Item(61, 9).Value = 1067.5  # I61: 1231.4546 -> 1067.5
$ws.Cells.Item(61, 10).Value = 2900  # J61: 1748 -> 2900
$ws.Cells.Item(61, 11).Value = 1067.5  # K61: 1231.4546 -> 1067.5
$ws.Cells.Item(61, 12).Value = 2900  # L61: 1748 -> 2900
$ws.Cells.Item(61, 13).Value = -855.5  # M61: -1019.4546 -> -855.5
$ws.Cells.Item(61, 14).Value = -3324  # N61: -2172 -> -3324
$ws.Cells.Item(74, 8).Value = 1608.3334  # H74: 842.75 -> 1608.3334
$ws.Cells.Item(74, 9).Value = 1608.3334  # I74: 793.3043 -> 1608.3334
$ws.Cells.Item(74, 10).Value = 0  # J74: 1980 -> 0
$ws.Cells.Item(74, 11).Value = 1608.3334  # K74: 793.3043 -> 1608.3334
$ws.Cells.Item(74, 12).Value = 0  # L74: 1980 -> 0
$ws.Cells.Item(74, 13).ClearContents()  # M74: remove (was 80.69569999999999)
$ws.Cells.Item(74, 14).Value = -734.3334  # N74: -3728 -> -734.3334
$ws.Cells.Item(77, 8).Value = 1608.3334  # H77: 842.75 -> 1608.3334
$ws.Cells.Item(77, 9).Value = 1608.3334  # I77: 793.3043 -> 1608.3334
$ws.Cells.Item(77, 10).Value = 0  # J77: 1980 -> 0
$ws.Cells.Item(77, 11).Value = 8041.666999999999  # K77: 3966.5215 -> 8041.666999999999
$ws.Cells.Item(77, 12).Value = 0  # L77: 9900 -> 0
$ws.Cells.Item(77, 13).ClearContents()  # M77: remove (was 401.4785000000002)
$ws.Cells.Item(77, 14).Value = -3673.666999999999  # N77: -18636 -> -3673.666999999999
$ws.Cells.Item(122, 8).Value = 2006  # H122: 759 -> 2006
$ws.Cells.Item(122, 9).Value = 2006  # I122: 703.8889 -> 2006
$ws.Cells.Item(122, 10).Value = 0  # J122: 1007 -> 0
$ws.Cells.Item(122, 11).Value = 6018  # K122: 2111.6667 -> 6018
$ws.Cells.Item(122, 12).Value = 0  # L122: 3021 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: remove (was 338.3332999999998)
$ws.Cells.Item(122, 14).Value = -3568  # N122: -7921 -> -3568
$ws.Cells.Item(136, 8).Value = 1271.1111  # H136: 1463.9 -> 1271.1111
$ws.Cells.Item(136, 9).Value = 1067.5  # I136: 1231.4546 -> 1067.5
$ws.Cells.Item(136, 10).Value = 2900  # J136: 1748 -> 2900
$ws.Cells.Item(136, 11).Value = 3202.5  # K136: 3694.3638 -> 3202.5
$ws.Cells.Item(136, 12).Value = 8700  # L136: 5244 -> 8700
$ws.Cells.Item(136, 13).Value = -652.5  # M136: -1144.3638 -> -652.5
$ws.Cells.Item(136, 14).Value = -13800  # N136: -10344 -> -13800

# --- Sheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 6579492.5  # H94: 6757315 -> 6579492.5
$ws.Cells.Item(94, 9).Value = 8929002  # I94: 9259703 -> 8929002
$ws.Cells.Item(94, 11).Value = 8929002  # K94: 9259703 -> 8929002
$ws.Cells.Item(94, 13).Value = -8928551  # M94: -9259252 -> -8928551
$ws.Cells.Item(134, 8).Value = 5222.2915  # H134: 5651.591 -> 5222.2915
$ws.Cells.Item(134, 9).Value = 1014.63635  # I134: 1066.1 -> 1014.63635
$ws.Cells.Item(134, 11).Value = 3043.90905  # K134: 3198.3 -> 3043.90905
$ws.Cells.Item(134, 13).Value = -508.9090500000002  # M134: -663.2999999999997 -> -508.9090500000002

# --- Sheet: CRP (63 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1121.4857  # H31: 1073.081 -> 1121.4857
$ws.Cells.Item(31, 9).Value = 763.89655  # I31: 729.19354 -> 763.89655
$ws.Cells.Item(31, 11).Value = 763.89655  # K31: 729.19354 -> 763.89655
$ws.Cells.Item(31, 13).Value = -468.89655  # M31: -434.19354 -> -468.89655
$ws.Cells.Item(32, 8).Value = 7950  # H32: 0 -> 7950
$ws.Cells.Item(32, 9).Value = 5900  # I32: 0 -> 5900
$ws.Cells.Item(32, 10).Value = 10000  # J32: 0 -> 10000
$ws.Cells.Item(32, 11).Value = 5900  # K32: 0 -> 5900
$ws.Cells.Item(32, 12).Value = 10000  # L32: 0 -> 10000
$ws.Cells.Item(32, 13).Value = -5584  # M32: None -> -5584
$ws.Cells.Item(32, 14).Value = -10632  # N32: None -> -10632
$ws.Cells.Item(34, 8).Value = 1121.4857  # H34: 1073.081 -> 1121.4857
$ws.Cells.Item(34, 9).Value = 763.89655  # I34: 729.19354 -> 763.89655
$ws.Cells.Item(34, 11).Value = 763.89655  # K34: 729.19354 -> 763.89655
$ws.Cells.Item(34, 13).Value = -561.89655  # M34: -527.19354 -> -561.89655
$ws.Cells.Item(52, 8).Value = 43623.168  # H52: 46989.75 -> 43623.168
$ws.Cells.Item(52, 10).Value = 43623.168  # J52: 46989.75 -> 43623.168
$ws.Cells.Item(52, 12).Value = 43623.168  # L52: 46989.75 -> 43623.168
$ws.Cells.Item(52, 14).Value = -44211.168  # N52: -47577.75 -> -44211.168
$ws.Cells.Item(76, 8).Value = 3693.3333  # H76: 0 -> 3693.3333
$ws.Cells.Item(76, 9).Value = 3693.3333  # I76: 0 -> 3693.3333
$ws.Cells.Item(76, 11).Value = 3693.3333  # K76: 0 -> 3693.3333
$ws.Cells.Item(76, 13).Value = -3378.3333  # M76: None -> -3378.3333
$ws.Cells.Item(79, 8).Value = 3693.3333  # H79: 0 -> 3693.3333
$ws.Cells.Item(79, 9).Value = 3693.3333  # I79: 0 -> 3693.3333
$ws.Cells.Item(79, 11).Value = 3693.3333  # K79: 0 -> 3693.3333
$ws.Cells.Item(79, 13).Value = -2601.3333  # M79: None -> -2601.3333
$ws.Cells.Item(99, 8).Value = 13158894  # H99: 2925771 -> 13158894
$ws.Cells.Item(99, 9).Value = 13158894  # I99: 4387654 -> 13158894
$ws.Cells.Item(99, 10).Value = 0  # J99: 2004.6666 -> 0
$ws.Cells.Item(99, 11).Value = 13158894  # K99: 4387654 -> 13158894
$ws.Cells.Item(99, 12).Value = 0  # L99: 2004.6666 -> 0
$ws.Cells.Item(99, 13).ClearContents()  # M99: remove (was -4386156)
$ws.Cells.Item(99, 14).Value = -13157396  # N99: -5000.6666 -> -13157396
$ws.Cells.Item(105, 8).Value = 966.5  # H105: 993.5 -> 966.5
$ws.Cells.Item(105, 9).Value = 966.5  # I105: 992.7143 -> 966.5
$ws.Cells.Item(105, 10).Value = 0  # J105: 999 -> 0
$ws.Cells.Item(105, 11).Value = 966.5  # K105: 992.7143 -> 966.5
$ws.Cells.Item(105, 12).Value = 0  # L105: 999 -> 0
$ws.Cells.Item(105, 13).ClearContents()  # M105: remove (was 754.2857)
$ws.Cells.Item(105, 14).Value = 780.5  # N105: -4493 -> 780.5
$ws.Cells.Item(122, 8).Value = 0  # H122: 845 -> 0
$ws.Cells.Item(122, 9).Value = 0  # I122: 845 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 2535 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: remove (was -85)
$ws.Cells.Item(126, 8).Value = 13158894  # H126: 2925771 -> 13158894
$ws.Cells.Item(126, 9).Value = 13158894  # I126: 4387654 -> 13158894
$ws.Cells.Item(126, 10).Value = 0  # J126: 2004.6666 -> 0
$ws.Cells.Item(126, 11).Value = 39476682  # K126: 13162962 -> 39476682
$ws.Cells.Item(126, 12).Value = 0  # L126: 6013.9998 -> 0
$ws.Cells.Item(126, 13).ClearContents()  # M126: remove (was -13160492)
$ws.Cells.Item(126, 14).Value = -39474212  # N126: -10953.9998 -> -39474212
$ws.Cells.Item(132, 8).Value = 10562.25  # H132: 12474.7 -> 10562.25
$ws.Cells.Item(132, 9).Value = 11774.8  # I132: 14468.5 -> 11774.8
$ws.Cells.Item(132, 11).Value = 35324.39999999999  # K132: 43405.5 -> 35324.39999999999
$ws.Cells.Item(132, 13).Value = -32794.39999999999  # M132: -40875.5 -> -32794.39999999999
$ws.Cells.Item(134, 8).Value = 14494418  # H134: 10102206 -> 14494418
$ws.Cells.Item(134, 9).Value = 19609378  # I134: 11905825 -> 19609378
$ws.Cells.Item(134, 10).Value = 2033.3334  # J134: 1940 -> 2033.3334
$ws.Cells.Item(134, 11).Value = 58828134  # K134: 35717475 -> 58828134
$ws.Cells.Item(134, 12).Value = 6100.0002  # L134: 5820 -> 6100.0002
$ws.Cells.Item(134, 13).Value = -58825599  # M134: -35714940 -> -58825599
$ws.Cells.Item(134, 14).Value = -11170.0002  # N134: -10890 -> -11170.0002

# --- Sheet: CUL (35 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 4546856  # H34: 3573101 -> 4546856
$ws.Cells.Item(34, 9).Value = 630  # I34: 662.5 -> 630
$ws.Cells.Item(34, 10).Value = 5883981  # J34: 4168507.5 -> 5883981
$ws.Cells.Item(34, 11).Value = 1890  # K34: 1987.5 -> 1890
$ws.Cells.Item(34, 12).Value = 17651943  # L34: 12505522.5 -> 17651943
$ws.Cells.Item(34, 13).Value = -1806  # M34: -1903.5 -> -1806
$ws.Cells.Item(34, 14).Value = -17652111  # N34: -12505690.5 -> -17652111
$ws.Cells.Item(39, 8).Value = 2668  # H39: 1000 -> 2668
$ws.Cells.Item(39, 10).Value = 3502  # J39: 0 -> 3502
$ws.Cells.Item(39, 12).Value = 10506  # L39: 0 -> 10506
$ws.Cells.Item(39, 14).Value = -11094  # N39: None -> -11094
$ws.Cells.Item(55, 8).Value = 3599.5  # H55: 3500 -> 3599.5
$ws.Cells.Item(55, 10).Value = 3599.5  # J55: 3500 -> 3599.5
$ws.Cells.Item(55, 12).Value = 10798.5  # L55: 10500 -> 10798.5
$ws.Cells.Item(55, 14).Value = -11152.5  # N55: -10854 -> -11152.5
$ws.Cells.Item(56, 8).Value = 7265.643  # H56: 7115.3076 -> 7265.643
$ws.Cells.Item(56, 9).Value = 7265.643  # I56: 7115.3076 -> 7265.643
$ws.Cells.Item(56, 11).Value = 7265.643  # K56: 7115.3076 -> 7265.643
$ws.Cells.Item(56, 13).Value = -6735.643  # M56: -6585.3076 -> -6735.643
$ws.Cells.Item(68, 8).Value = 2057.6304  # H68: 2110.2727 -> 2057.6304
$ws.Cells.Item(68, 10).Value = 2092.2444  # J68: 2147.721 -> 2092.2444
$ws.Cells.Item(68, 12).Value = 6276.733200000001  # L68: 6443.163 -> 6276.733200000001
$ws.Cells.Item(68, 14).Value = -7898.733200000001  # N68: -8065.163 -> -7898.733200000001
$ws.Cells.Item(71, 8).Value = 2057.6304  # H71: 2110.2727 -> 2057.6304
$ws.Cells.Item(71, 10).Value = 2092.2444  # J71: 2147.721 -> 2092.2444
$ws.Cells.Item(71, 12).Value = 18830.1996  # L71: 19329.489 -> 18830.1996
$ws.Cells.Item(71, 14).Value = -26942.1996  # N71: -27441.489 -> -26942.1996
$ws.Cells.Item(137, 8).Value = 8995.08  # H137: 9286.541999999999 -> 8995.08
$ws.Cells.Item(137, 10).Value = 13043.5625  # J137: 13779.8 -> 13043.5625
$ws.Cells.Item(137, 12).Value = 39130.6875  # L137: 41339.39999999999 -> 39130.6875
$ws.Cells.Item(137, 14).Value = -49330.6875  # N137: -51539.39999999999 -> -49330.6875
$ws.Cells.Item(141, 8).Value = 1726.6  # H141: 1781.4546 -> 1726.6
$ws.Cells.Item(141, 9).Value = 1726.6  # I141: 1781.4546 -> 1726.6
$ws.Cells.Item(141, 11).Value = 5179.799999999999  # K141: 5344.3638 -> 5179.799999999999
$ws.Cells.Item(141, 13).Value = 0.2000000000007276  # M141: -164.3638000000001 -> 0.2000000000007276

# --- Sheet: GSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 29999.8  # H93: 0 -> 29999.8
$ws.Cells.Item(93, 10).Value = 29999.8  # J93: 0 -> 29999.8
$ws.Cells.Item(93, 12).Value = 29999.8  # L93: 0 -> 29999.8
$ws.Cells.Item(93, 14).Value = -33743.8  # N93: None -> -33743.8
$ws.Cells.Item(104, 8).Value = 50667.75  # H104: 53800 -> 50667.75
$ws.Cells.Item(104, 10).Value = 50667.75  # J104: 53800 -> 50667.75
$ws.Cells.Item(104, 12).Value = 50667.75  # L104: 53800 -> 50667.75
$ws.Cells.Item(104, 14).Value = -57655.75  # N104: -60788 -> -57655.75
$ws.Cells.Item(113, 8).Value = 2499.2856  # H113: 2328.3333 -> 2499.2856
$ws.Cells.Item(113, 9).Value = 1874  # I113: 1770.5 -> 1874
$ws.Cells.Item(113, 10).Value = 2968.25  # J113: 2774.6 -> 2968.25
$ws.Cells.Item(113, 11).Value = 1874  # K113: 1770.5 -> 1874
$ws.Cells.Item(113, 12).Value = 2968.25  # L113: 2774.6 -> 2968.25
$ws.Cells.Item(113, 13).Value = 296  # M113: 399.5 -> 296
$ws.Cells.Item(113, 14).Value = -7308.25  # N113: -7114.6 -> -7308.25
$ws.Cells.Item(132, 8).Value = 2975.4443  # H132: 2335.606 -> 2975.4443
$ws.Cells.Item(132, 9).Value = 2472.5  # I132: 1744.3334 -> 2472.5
$ws.Cells.Item(132, 10).Value = 6999  # J132: 4996.3335 -> 6999
$ws.Cells.Item(132, 11).Value = 7417.5  # K132: 5233.0002 -> 7417.5
$ws.Cells.Item(132, 12).Value = 20997  # L132: 14989.0005 -> 20997
$ws.Cells.Item(132, 13).Value = -4887.5  # M132: -2703.0002 -> -4887.5
$ws.Cells.Item(132, 14).Value = -26057  # N132: -20049.0005 -> -26057

# --- Sheet: LTW (8 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2000.2222  # H16: 2072.2222 -> 2000.2222
$ws.Cells.Item(16, 9).Value = 1937.75  # I16: 2018.75 -> 1937.75
$ws.Cells.Item(16, 11).Value = 1937.75  # K16: 2018.75 -> 1937.75
$ws.Cells.Item(16, 13).Value = -1767.75  # M16: -1848.75 -> -1767.75
$ws.Cells.Item(46, 8).Value = 6621.5  # H46: 7076.923 -> 6621.5
$ws.Cells.Item(46, 9).Value = 725.25  # I46: 733.3333 -> 725.25
$ws.Cells.Item(46, 11).Value = 725.25  # K46: 733.3333 -> 725.25
$ws.Cells.Item(46, 13).Value = -537.25  # M46: -545.3333 -> -537.25

# --- Sheet: WVR (18 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2256.125  # H132: 2413.75 -> 2256.125
$ws.Cells.Item(132, 9).Value = 2124.4482  # I132: 2296.5 -> 2124.4482
$ws.Cells.Item(132, 10).Value = 2603.2727  # J132: 2718.6 -> 2603.2727
$ws.Cells.Item(132, 11).Value = 6373.344599999999  # K132: 6889.5 -> 6373.344599999999
$ws.Cells.Item(132, 12).Value = 7809.8181  # L132: 8155.799999999999 -> 7809.8181
$ws.Cells.Item(132, 13).Value = -3843.344599999999  # M132: -4359.5 -> -3843.344599999999
$ws.Cells.Item(132, 14).Value = -12869.8181  # N132: -13215.8 -> -12869.8181
$ws.Cells.Item(136, 8).Value = 556.575  # H136: 484.02942 -> 556.575
$ws.Cells.Item(136, 9).Value = 300.5  # I136: 307.96155 -> 300.5
$ws.Cells.Item(136, 10).Value = 1032.1428  # J136: 1056.25 -> 1032.1428
$ws.Cells.Item(136, 11).Value = 901.5  # K136: 923.88465 -> 901.5
$ws.Cells.Item(136, 12).Value = 3096.4284  # L136: 3168.75 -> 3096.4284
$ws.Cells.Item(136, 13).Value = 1648.5  # M136: 1626.11535 -> 1648.5
$ws.Cells.Item(136, 14).Value = -8196.428400000001  # N136: -8268.75 -> -8196.428400000001
$ws.Cells.Item(139, 8).Value = 49015  # H139: 0 -> 49015
$ws.Cells.Item(139, 10).Value = 49015  # J139: 0 -> 49015
$ws.Cells.Item(139, 12).Value = 49015  # L139: 0 -> 49015
$ws.Cells.Item(139, 14).Value = -59295  # N139: None -> -59295
